$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 901.25
$ws.Range("I107").Value = 1332.8
$ws.Range("J107").Value = 469.7
$ws.Range("K107").Value = 1332.8
$ws.Range("L107").Value = 469.7
$ws.Range("M107").Value = 587.2
$ws.Range("N107").Value = -4309.7
$ws.Range("H116").Value = 3580.8667
$ws.Range("I116").Value = 3843
$ws.Range("J116").Value = 3281.2856
$ws.Range("K116").Value = 3843
$ws.Range("L116").Value = 3281.2856
$ws.Range("M116").Value = -401
$ws.Range("N116").Value = -10165.2856
$ws.Range("H132").Value = 5661.175
$ws.Range("I132").Value = 2124.9
$ws.Range("J132").Value = 16270
$ws.Range("K132").Value = 6374.700000000001
$ws.Range("L132").Value = 48810
$ws.Range("M132").Value = -3844.700000000001
$ws.Range("N132").Value = -53870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1640.8182
$ws.Range("I2").Value = 1890.375
$ws.Range("K2").Value = 1890.375
$ws.Range("M2").Value = -1777.375
$ws.Range("H45").Value = 669.1177
$ws.Range("I45").Value = 627.7692
$ws.Range("J45").Value = 803.5
$ws.Range("K45").Value = 627.7692
$ws.Range("L45").Value = 803.5
$ws.Range("M45").Value = -250.7692
$ws.Range("N45").Value = -1557.5
$ws.Range("H74").Value = 11110561
$ws.Range("I74").Value = 8804088
$ws.Range("J74").Value = 15199309
$ws.Range("K74").Value = 8804088
$ws.Range("L74").Value = 15199309
$ws.Range("M74").Value = -8803214
$ws.Range("N74").Value = -15201057
$ws.Range("H77").Value = 11110561
$ws.Range("I77").Value = 8804088
$ws.Range("J77").Value = 15199309
$ws.Range("K77").Value = 44020440
$ws.Range("L77").Value = 75996545
$ws.Range("M77").Value = -44016072
$ws.Range("N77").Value = -76005281
$ws.Range("H116").Value = 1640.8182
$ws.Range("I116").Value = 1890.375
$ws.Range("K116").Value = 1890.375
$ws.Range("M116").Value = 403.625
$ws.Range("H132").Value = 28492.205
$ws.Range("I132").Value = 61335.47
$ws.Range("J132").Value = 3113.318
$ws.Range("K132").Value = 184006.41
$ws.Range("L132").Value = 9339.954000000002
$ws.Range("M132").Value = -181476.41
$ws.Range("N132").Value = -14399.954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1640.8182
$ws.Range("I3").Value = 1890.375
$ws.Range("K3").Value = 1890.375
$ws.Range("M3").Value = -1776.375
$ws.Range("H80").Value = 236.36363
$ws.Range("I80").Value = 297.85715
$ws.Range("J80").Value = 128.75
$ws.Range("K80").Value = 297.85715
$ws.Range("L80").Value = 128.75
$ws.Range("M80").Value = 700.14285
$ws.Range("N80").Value = -2124.75
$ws.Range("H83").Value = 236.36363
$ws.Range("I83").Value = 297.85715
$ws.Range("J83").Value = 128.75
$ws.Range("K83").Value = 1489.28575
$ws.Range("L83").Value = 643.75
$ws.Range("M83").Value = 3502.71425
$ws.Range("N83").Value = -10627.75
$ws.Range("H99").Value = 15401794
$ws.Range("I99").Value = 5516506.5
$ws.Range("J99").Value = 50000300
$ws.Range("K99").Value = 5516506.5
$ws.Range("L99").Value = 50000300
$ws.Range("M99").Value = -5515008.5
$ws.Range("N99").Value = -50003296
$ws.Range("H105").Value = 1875.4231
$ws.Range("I105").Value = 1780.909
$ws.Range("J105").Value = 2395.25
$ws.Range("K105").Value = 1780.909
$ws.Range("L105").Value = 2395.25
$ws.Range("M105").Value = -33.90900000000011
$ws.Range("N105").Value = -5889.25
$ws.Range("H107").Value = 749.2105
$ws.Range("I107").Value = 763.05554
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 763.05554
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1156.94446
$ws.Range("N107").Value = -4340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1400.1818
$ws.Range("I22").Value = 1450.2
$ws.Range("K22").Value = 1450.2
$ws.Range("M22").Value = -1100.2
$ws.Range("H31").Value = 1923.29
$ws.Range("I31").Value = 726.61664
$ws.Range("J31").Value = 3718.3
$ws.Range("K31").Value = 726.61664
$ws.Range("L31").Value = 3718.3
$ws.Range("M31").Value = -431.61664
$ws.Range("N31").Value = -4308.3
$ws.Range("H34").Value = 1923.29
$ws.Range("I34").Value = 726.61664
$ws.Range("J34").Value = 3718.3
$ws.Range("K34").Value = 726.61664
$ws.Range("L34").Value = 3718.3
$ws.Range("M34").Value = -524.61664
$ws.Range("N34").Value = -4122.3
$ws.Range("H86").Value = 4868.909
$ws.Range("I86").Value = 3224.2727
$ws.Range("J86").Value = 6513.5454
$ws.Range("K86").Value = 3224.2727
$ws.Range("L86").Value = 6513.5454
$ws.Range("M86").Value = -2101.2727
$ws.Range("N86").Value = -8759.545399999999
$ws.Range("H89").Value = 4868.909
$ws.Range("I89").Value = 3224.2727
$ws.Range("J89").Value = 6513.5454
$ws.Range("K89").Value = 16121.3635
$ws.Range("L89").Value = 32567.727
$ws.Range("M89").Value = -10505.3635
$ws.Range("N89").Value = -43799.727
$ws.Range("H105").Value = 1149.1389
$ws.Range("I105").Value = 992.875
$ws.Range("J105").Value = 1461.6666
$ws.Range("K105").Value = 992.875
$ws.Range("L105").Value = 1461.6666
$ws.Range("M105").Value = 754.125
$ws.Range("N105").Value = -4955.6666
$ws.Range("H107").Value = 302.6857
$ws.Range("I107").Value = 291.7619
$ws.Range("J107").Value = 319.07144
$ws.Range("K107").Value = 291.7619
$ws.Range("L107").Value = 319.07144
$ws.Range("M107").Value = 1628.2381
$ws.Range("N107").Value = -4159.07144
$ws.Range("H122").Value = 1675.5883
$ws.Range("I122").Value = 1830.4667
$ws.Range("J122").Value = 514
$ws.Range("K122").Value = 5491.4001
$ws.Range("L122").Value = 1542
$ws.Range("M122").Value = -3041.4001
$ws.Range("N122").Value = -6442
$ws.Range("H132").Value = 14708658
$ws.Range("I132").Value = 20835494
$ws.Range("J132").Value = 4251.2
$ws.Range("K132").Value = 62506482
$ws.Range("L132").Value = 12753.6
$ws.Range("M132").Value = -62503952
$ws.Range("N132").Value = -17813.6
$ws.Range("H134").Value = 17243658
$ws.Range("I134").Value = 26317560
$ws.Range("J134").Value = 3240.2
$ws.Range("K134").Value = 78952680
$ws.Range("L134").Value = 9720.599999999999
$ws.Range("M134").Value = -78950145
$ws.Range("N134").Value = -14790.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1316
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 1297.1282
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 3891.3846
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -5513.3846
$ws.Range("H71").Value = 1316
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 1297.1282
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 11674.1538
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -19786.1538
$ws.Range("H80").Value = 2133.3333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2133.3333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6399.999899999999
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -8271.999899999999
$ws.Range("H83").Value = 2133.3333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2133.3333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 19199.9997
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -28559.9997
$ws.Range("H92").Value = 333.83334
$ws.Range("J92").Value = 325.75
$ws.Range("L92").Value = 977.25
$ws.Range("N92").Value = -3473.25
$ws.Range("H114").Value = 373.42856
$ws.Range("I114").Value = 373.42856
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 1120.28568
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 2133.71432
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1210.909
$ws.Range("I122").Value = 1213.3334
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3640.0002
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1190.0002
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21580868
$ws.Range("I136").Value = 33705180
$ws.Range("J136").Value = 5047718
$ws.Range("K136").Value = 101115540
$ws.Range("L136").Value = 15143154
$ws.Range("M136").Value = -101112990
$ws.Range("N136").Value = -15148254
